# Reorder the "manifest" columns on the "all" sheet: the column that used to
# sit in F ("assumed_answer" -> renamed "assumed_easy") and the column that
# used to sit in G ("col_id_easy") trade places, F becomes col_id_easy and
# G becomes assumed_easy (widths travel with the content).
#
# Using Cut/Insert (the classic Excel "swap columns" gesture) moves column F
# to sit just before the new H, which pushes the old G left into F and lands
# the old F's data in G - i.e. a clean F<->G swap, matching the target layout:
#   F1: assumed_answer -> col_id_easy
#   G1: col_id_easy    -> assumed_answer (renamed below to assumed_easy)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all")

$ws.Columns.Item(6).Cut() | Out-Null
$ws.Columns.Item(8).Insert() | Out-Null

# The moved header text is also being renamed as part of this edit.
$ws.Range("G1").Value = "assumed_easy"
